# Revision gestion des risques
#
# Renames the "Autocorrélation " sheet to "Accordeur" and updates the
# selection/active-cell state on that sheet to match the new view
# (A6:M12 selected, A6 active) as captured in the target workbook.

$wb = $excel.ActiveWorkbook

# --- Locate the "Autocorrélation " sheet (fall back to the 2nd sheet) ---
$ws2 = $null
foreach ($s in $wb.Worksheets) {
    if ($s.Name -eq "Autocorrélation ") {
        $ws2 = $s
    }
}
if ($ws2 -eq $null) {
    $ws2 = $wb.Worksheets.Item(2)
}

# --- Rename it to "Accordeur" ---
$ws2.Name = "Accordeur"

# --- Update the active selection on that sheet ---
$ws2.Activate()
[void]$ws2.Range("A6:M12").Select()
